$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Helper: split the current Find match off from its neighbouring text
# by toggling a character property on/off. Word always has to end the
# run at the edges of a reformatted sub-range, so flipping Bold true
# then back to false on the exact sub-range forces a run break there
# without leaving any residual formatting difference (Bold False is
# the default and is omitted from the saved rPr).
# ------------------------------------------------------------------

# 1) Split every "copy+paste" occurrence away from its neighbours, so
#    it becomes its own run (this mirrors the proofErr-wrapped run the
#    diff introduces around the English word "copy+paste").
$rng = $d.Content
while ($rng.Find.Execute("copy+paste", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Bold = 1
    $rng.Bold = 0
    $rng.Collapse(0)
}

# 2) Split "triggers" and "cursors" out of "Υλοποίηση triggers, cursors"
#    into their own runs too.
$rng = $d.Content
while ($rng.Find.Execute("triggers, cursors", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $matchStart = $rng.Start
    $matchEnd = $rng.End

    $triggersRng = $d.Range($matchStart, $matchStart + 8)
    $triggersRng.Bold = 1
    $triggersRng.Bold = 0

    $cursorsRng = $d.Range($matchStart + 10, $matchEnd)
    $cursorsRng.Bold = 1
    $cursorsRng.Bold = 0

    $rng = $d.Range($matchEnd, $d.Content.End)
}

# 3) pgadmin / sql / db / py runs only gain proofErr spell-check
#    markers in the diff (no visible text changes) -- nothing further
#    to do for them here.

# 4) "db_api.py" -> "db_client.py": the lone lowercase run "api"
#    becomes "client".
$d.Content.Find.Execute("api", $true, $false, $false, $false, $false, $true, 1, $false, "client", 2) | Out-Null

Write-Output "edit complete"
